# Auto-generated edit: realigns data rows 2-16 of the active worksheet so that each
# "Beteckning" record (columns A-Z) matches the refreshed scrape snapshot. The same
# 15 case records are kept, only their row order changes (plus the shared "Forandrad"
# date in column C advances by one day, and the extra column-Z hyperlink formula
# (Fageltillsynsbegaranslank) now belongs to the row holding "A 13467-2023").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target row 2  <=  data for "A 13467-2023" (previously row 3)
$ws.Range("A2").Value = "A 13467-2023"
$ws.Range("B2").Value = 45005
$ws.Range("C2").Value = 46065
$ws.Range("D2").Value = "SKÅNE LÄN"
$ws.Range("E2").Value = "HELSINGBORG"
$ws.Range("G2").Value = 2.3
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = "Gulsparv`r`nHypoxylon petriniae`r`nKråka`r`nGrå skärelav`r`nGulnål"
$ws.Range("S2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/artfynd/A 13467-2023 artfynd.xlsx"", ""A 13467-2023"")"
$ws.Range("T2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/kartor/A 13467-2023 karta.png"", ""A 13467-2023"")"
$ws.Range("V2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/klagomål/A 13467-2023 FSC-klagomål.docx"", ""A 13467-2023"")"
$ws.Range("W2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/klagomålsmail/A 13467-2023 FSC-klagomål mail.docx"", ""A 13467-2023"")"
$ws.Range("X2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/tillsyn/A 13467-2023 tillsynsbegäran.docx"", ""A 13467-2023"")"
$ws.Range("Y2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/tillsynsmail/A 13467-2023 tillsynsbegäran mail.docx"", ""A 13467-2023"")"
$ws.Range("Z2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/fåglar/A 13467-2023 prioriterade fågelarter.docx"", ""A 13467-2023"")"

# Target row 3  <=  data for "A 45325-2025" (previously row 2)
$ws.Range("A3").Value = "A 45325-2025"
$ws.Range("B3").Value = 45922
$ws.Range("C3").Value = 46065
$ws.Range("D3").Value = "SKÅNE LÄN"
$ws.Range("E3").Value = "HELSINGBORG"
$ws.Range("G3").Value = 1.6
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = "Nordlig buksimmare`r`nStörre vattensalamander`r`nÅkergroda`r`nMindre vattensalamander`r`nVanlig groda"
$ws.Range("S3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/artfynd/A 45325-2025 artfynd.xlsx"", ""A 45325-2025"")"
$ws.Range("T3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/kartor/A 45325-2025 karta.png"", ""A 45325-2025"")"
$ws.Range("V3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/klagomål/A 45325-2025 FSC-klagomål.docx"", ""A 45325-2025"")"
$ws.Range("W3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/klagomålsmail/A 45325-2025 FSC-klagomål mail.docx"", ""A 45325-2025"")"
$ws.Range("X3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/tillsyn/A 45325-2025 tillsynsbegäran.docx"", ""A 45325-2025"")"
$ws.Range("Y3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/tillsynsmail/A 45325-2025 tillsynsbegäran mail.docx"", ""A 45325-2025"")"
$ws.Range("Z3").Formula = ""

# Target row 4  <=  data for "A 45832-2023" (previously row 4)
$ws.Range("A4").Value = "A 45832-2023"
$ws.Range("B4").Value = 45195
$ws.Range("C4").Value = 46065
$ws.Range("D4").Value = "SKÅNE LÄN"
$ws.Range("E4").Value = "HELSINGBORG"
$ws.Range("G4").Value = 2.3
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Skånebjörnbär"
$ws.Range("S4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/artfynd/A 45832-2023 artfynd.xlsx"", ""A 45832-2023"")"
$ws.Range("T4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/kartor/A 45832-2023 karta.png"", ""A 45832-2023"")"
$ws.Range("V4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/klagomål/A 45832-2023 FSC-klagomål.docx"", ""A 45832-2023"")"
$ws.Range("W4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/klagomålsmail/A 45832-2023 FSC-klagomål mail.docx"", ""A 45832-2023"")"
$ws.Range("X4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/tillsyn/A 45832-2023 tillsynsbegäran.docx"", ""A 45832-2023"")"
$ws.Range("Y4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1283/tillsynsmail/A 45832-2023 tillsynsbegäran mail.docx"", ""A 45832-2023"")"
$ws.Range("Z4").Formula = ""

# Target row 5  <=  data for "A 24227-2022" (previously row 5)
$ws.Range("A5").Value = "A 24227-2022"
$ws.Range("B5").Value = 44725.64246527778
$ws.Range("C5").Value = 46065
$ws.Range("D5").Value = "SKÅNE LÄN"
$ws.Range("E5").Value = "HELSINGBORG"
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").ClearContents()
$ws.Range("S5").ClearContents()
$ws.Range("T5").ClearContents()
$ws.Range("V5").ClearContents()
$ws.Range("W5").ClearContents()
$ws.Range("X5").ClearContents()
$ws.Range("Y5").ClearContents()
$ws.Range("Z5").Formula = ""

# Target row 6  <=  data for "A 5792-2024" (previously row 8)
$ws.Range("A6").Value = "A 5792-2024"
$ws.Range("B6").Value = 45335
$ws.Range("C6").Value = 46065
$ws.Range("D6").Value = "SKÅNE LÄN"
$ws.Range("E6").Value = "HELSINGBORG"
$ws.Range("G6").Value = 5.6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").ClearContents()
$ws.Range("S6").ClearContents()
$ws.Range("T6").ClearContents()
$ws.Range("V6").ClearContents()
$ws.Range("W6").ClearContents()
$ws.Range("X6").ClearContents()
$ws.Range("Y6").ClearContents()
$ws.Range("Z6").Formula = ""

# Target row 7  <=  data for "A 13651-2023" (previously row 9)
$ws.Range("A7").Value = "A 13651-2023"
$ws.Range("B7").Value = 45006
$ws.Range("C7").Value = 46065
$ws.Range("D7").Value = "SKÅNE LÄN"
$ws.Range("E7").Value = "HELSINGBORG"
$ws.Range("G7").Value = 2.2
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("V7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").Formula = ""

# Target row 8  <=  data for "A 12651-2022" (previously row 7)
$ws.Range("A8").Value = "A 12651-2022"
$ws.Range("B8").Value = 44641
$ws.Range("C8").Value = 46065
$ws.Range("D8").Value = "SKÅNE LÄN"
$ws.Range("E8").Value = "HELSINGBORG"
$ws.Range("G8").Value = 3.2
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("V8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").Formula = ""

# Target row 9  <=  data for "A 2593-2024" (previously row 6)
$ws.Range("A9").Value = "A 2593-2024"
$ws.Range("B9").Value = 45313.69204861111
$ws.Range("C9").Value = 46065
$ws.Range("D9").Value = "SKÅNE LÄN"
$ws.Range("E9").Value = "HELSINGBORG"
$ws.Range("G9").Value = 2.3
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("V9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").Formula = ""

# Target row 10  <=  data for "A 7333-2025" (previously row 16)
$ws.Range("A10").Value = "A 7333-2025"
$ws.Range("B10").Value = 45703.35899305555
$ws.Range("C10").Value = 46065
$ws.Range("D10").Value = "SKÅNE LÄN"
$ws.Range("E10").Value = "HELSINGBORG"
$ws.Range("G10").Value = 0.9
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").ClearContents()
$ws.Range("S10").ClearContents()
$ws.Range("T10").ClearContents()
$ws.Range("V10").ClearContents()
$ws.Range("W10").ClearContents()
$ws.Range("X10").ClearContents()
$ws.Range("Y10").ClearContents()
$ws.Range("Z10").Formula = ""

# Target row 11  <=  data for "A 35642-2023" (previously row 12)
$ws.Range("A11").Value = "A 35642-2023"
$ws.Range("B11").Value = 45147
$ws.Range("C11").Value = 46065
$ws.Range("D11").Value = "SKÅNE LÄN"
$ws.Range("E11").Value = "HELSINGBORG"
$ws.Range("G11").Value = 1.2
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("R11").ClearContents()
$ws.Range("S11").ClearContents()
$ws.Range("T11").ClearContents()
$ws.Range("V11").ClearContents()
$ws.Range("W11").ClearContents()
$ws.Range("X11").ClearContents()
$ws.Range("Y11").ClearContents()
$ws.Range("Z11").Formula = ""

# Target row 12  <=  data for "A 28288-2023" (previously row 15)
$ws.Range("A12").Value = "A 28288-2023"
$ws.Range("B12").Value = 45099.6349537037
$ws.Range("C12").Value = 46065
$ws.Range("D12").Value = "SKÅNE LÄN"
$ws.Range("E12").Value = "HELSINGBORG"
$ws.Range("G12").Value = 0.5
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").ClearContents()
$ws.Range("S12").ClearContents()
$ws.Range("T12").ClearContents()
$ws.Range("V12").ClearContents()
$ws.Range("W12").ClearContents()
$ws.Range("X12").ClearContents()
$ws.Range("Y12").ClearContents()
$ws.Range("Z12").Formula = ""

# Target row 13  <=  data for "A 8194-2025" (previously row 10)
$ws.Range("A13").Value = "A 8194-2025"
$ws.Range("B13").Value = 45708
$ws.Range("C13").Value = 46065
$ws.Range("D13").Value = "SKÅNE LÄN"
$ws.Range("E13").Value = "HELSINGBORG"
$ws.Range("G13").Value = 1.9
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").ClearContents()
$ws.Range("S13").ClearContents()
$ws.Range("T13").ClearContents()
$ws.Range("V13").ClearContents()
$ws.Range("W13").ClearContents()
$ws.Range("X13").ClearContents()
$ws.Range("Y13").ClearContents()
$ws.Range("Z13").Formula = ""

# Target row 14  <=  data for "A 50997-2025" (previously row 11)
$ws.Range("A14").Value = "A 50997-2025"
$ws.Range("B14").Value = 45946
$ws.Range("C14").Value = 46065
$ws.Range("D14").Value = "SKÅNE LÄN"
$ws.Range("E14").Value = "HELSINGBORG"
$ws.Range("G14").Value = 1.5
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("R14").ClearContents()
$ws.Range("S14").ClearContents()
$ws.Range("T14").ClearContents()
$ws.Range("V14").ClearContents()
$ws.Range("W14").ClearContents()
$ws.Range("X14").ClearContents()
$ws.Range("Y14").ClearContents()
$ws.Range("Z14").Formula = ""

# Target row 15  <=  data for "A 7814-2026" (previously row 14)
$ws.Range("A15").Value = "A 7814-2026"
$ws.Range("B15").Value = 46062.61388888889
$ws.Range("C15").Value = 46065
$ws.Range("D15").Value = "SKÅNE LÄN"
$ws.Range("E15").Value = "HELSINGBORG"
$ws.Range("G15").Value = 1.1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("R15").ClearContents()
$ws.Range("S15").ClearContents()
$ws.Range("T15").ClearContents()
$ws.Range("V15").ClearContents()
$ws.Range("W15").ClearContents()
$ws.Range("X15").ClearContents()
$ws.Range("Y15").ClearContents()
$ws.Range("Z15").Formula = ""

# Target row 16  <=  data for "A 7827-2026" (previously row 13)
$ws.Range("A16").Value = "A 7827-2026"
$ws.Range("B16").Value = 46062.63958333333
$ws.Range("C16").Value = 46065
$ws.Range("D16").Value = "SKÅNE LÄN"
$ws.Range("E16").Value = "HELSINGBORG"
$ws.Range("G16").Value = 2.1
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("R16").ClearContents()
$ws.Range("S16").ClearContents()
$ws.Range("T16").ClearContents()
$ws.Range("V16").ClearContents()
$ws.Range("W16").ClearContents()
$ws.Range("X16").ClearContents()
$ws.Range("Y16").ClearContents()
$ws.Range("Z16").Formula = ""

